# Weekly update: two new price records for "Alcachofa" at Vega Modelo de
# Temuco were reported. Insert them as new rows 190-191 (pushing the
# existing rows 190-204 down to 192-206, preserving all their data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("190:191").Insert()

# New row 190
$ws.Range("A190").Value = 10
$ws.Range("B190").Value = 'Vega Modelo de Temuco'
$ws.Range("C190").Value = 'La Araucanía'
$ws.Range("D190").Value = 44783
$ws.Range("E190").Value = 9
$ws.Range("F190").Value = 100112013
$ws.Range("G190").Value = 'Alcachofa'
$ws.Range("H190").Value = 'Española'
$ws.Range("I190").Value = 'Primera'
$ws.Range("J190").Value = 120
$ws.Range("K190").Value = 18000
$ws.Range("L190").Value = 18000
$ws.Range("M190").Value = 18000
$ws.Range("N190").Value = '$/caja 30 unidades'
$ws.Range("O190").Value = 'Provincia de Limarí'
$ws.Range("P190").Value = 600
$ws.Range("Q190").Value = 30
$ws.Range("R190").Value = 'Hortaliza'

# New row 191
$ws.Range("A191").Value = 10
$ws.Range("B191").Value = 'Vega Modelo de Temuco'
$ws.Range("C191").Value = 'La Araucanía'
$ws.Range("D191").Value = 44783
$ws.Range("E191").Value = 9
$ws.Range("F191").Value = 100112013
$ws.Range("G191").Value = 'Alcachofa'
$ws.Range("H191").Value = 'Madrigal'
$ws.Range("I191").Value = 'Primera'
$ws.Range("J191").Value = 200
$ws.Range("K191").Value = 15000
$ws.Range("L191").Value = 15000
$ws.Range("M191").Value = 15000
$ws.Range("N191").Value = '$/caja 40 unidades'
$ws.Range("O191").Value = 'Provincia de Limarí'
$ws.Range("P191").Value = 375
$ws.Range("Q191").Value = 40
$ws.Range("R191").Value = 'Hortaliza'
